# Update the crypto-price snapshot (Price / Volume(1h) columns) with the
# latest scrape values, per "Updated symbol list ... with GitHub Actions".
#
# NOTE: columns D (Price) and E (Volume(1h)) are stored as literal TEXT in
# this sheet (not numbers), so values such as "261.44" or "0.06%" must stay
# text. Plain `Range.Value = "261.44"` would be auto-coerced to a number
# (and "0.06%" to a fraction), losing the exact display string, so each
# literal is written with a leading apostrophe to force text entry - the
# same trick Excel itself uses - which keeps the stored text identical to
# the source values while marking the cell as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.44"
$ws.Range("E2").Value = "'0.06%"

$ws.Range("D3").Value = "'26.75"
$ws.Range("E3").Value = "'-2.11%"

$ws.Range("D4").Value = "'4.707"
$ws.Range("E4").Value = "'-0.07%"

$ws.Range("D5").Value = "'0.06198"
$ws.Range("E5").Value = "'1.75%"

$ws.Range("D6").Value = "'6.748"
$ws.Range("E6").Value = "'0.99%"

$ws.Range("D7").Value = "'0.8487"
$ws.Range("E7").Value = "'0.08%"

$ws.Range("E8").Value = "'-1.85%"

$ws.Range("D9").Value = "'0.1405"
$ws.Range("E9").Value = "'-0.16%"

$ws.Range("D10").Value = "'0.05090"
$ws.Range("E10").Value = "'6.83%"

$ws.Range("D11").Value = "'0.07105"
$ws.Range("E11").Value = "'0.07%"

$ws.Range("D12").Value = "'0.03090"
$ws.Range("E12").Value = "'-0.04%"

$ws.Range("D13").Value = "'0.09043"
$ws.Range("E13").Value = "'-0.24%"

$ws.Range("D14").Value = "'0.001533"
$ws.Range("E14").Value = "'-0.61%"

$ws.Range("D15").Value = "'0.0006163"
$ws.Range("E15").Value = "'1.08%"

$ws.Range("D16").Value = "'0.005943"
$ws.Range("E16").Value = "'-3.27%"

$ws.Range("D17").Value = "'3.448"
$ws.Range("E17").Value = "'0.00%"

$ws.Range("D18").Value = "'3.175"
$ws.Range("E18").Value = "'0.90%"

$ws.Range("D19").Value = "'2.187"
$ws.Range("E19").Value = "'1.12%"

$ws.Range("E20").Value = "'-0.21%"

$ws.Range("E21").Value = "'0.41%"

$ws.Range("D22").Value = "'4.106"
$ws.Range("E22").Value = "'0.27%"

$ws.Range("D23").Value = "'0.04236"
$ws.Range("E23").Value = "'0.08%"

$ws.Range("D24").Value = "'0.001177"
$ws.Range("E24").Value = "'-3.75%"

$ws.Range("E25").Value = "'6.96%"

$ws.Range("E26").Value = "'0.05%"

$ws.Range("D40").Value = "'0.03956"
$ws.Range("E40").Value = "'2.04%"

$ws.Range("D41").Value = "'0.1112"
$ws.Range("E41").Value = "'-0.08%"

$ws.Range("D42").Value = "'0.004134"
$ws.Range("E42").Value = "'1.20%"

$ws.Range("D43").Value = "'0.002133"
$ws.Range("E43").Value = "'-3.77%"

$ws.Range("D44").Value = "'0.01327"
$ws.Range("E44").Value = "'-18.41%"

$ws.Range("D45").Value = "'0.00005165"
$ws.Range("E45").Value = "'0.54%"

$ws.Range("E46").Value = "'0.05%"

$ws.Range("D48").Value = "'0.2497"
$ws.Range("E48").Value = "'82.09%"

$ws.Range("E49").Value = "'0.05%"

$ws.Range("E50").Value = "'0.05%"
